$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank placeholder row 2, shifting all data rows up by one.
$ws.Rows.Item(2).Delete()

# Freeze panes below the header row (row 1), matching the new sheetView pane/selection.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("B3").Select()
